$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "X1"
$ws.Range("C1").Value = "X2"
$ws.Range("D1").Value = "X3"
$ws.Range("E1").Value = "X4"

# ---- Data rows ----
$ws.Range("A2").Value = "Impacto1"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.7
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = "Impacto2"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 1

# Remove the old row 4 values entirely (A4:C4 previously held 1, 0, -15)
$ws.Range("A4:C4").ClearContents()

# ---- Formatting ----
# Header font: size 10, black (built on a single cell first, then propagated via
# a format-only paste so we don't leave behind extra unused cell-format records)
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Font.Color = 0
$ws.Range("B1").Copy()
$ws.Range("C1:AK1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Number format 0.0 across the data block (rows 2-4, columns B-AK)
$ws.Range("B2:AK4").NumberFormat = "0.0"

# ---- Selection ----
[void]$ws.Range("F3").Select()
